$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 142.8125
$ws.Range("I33").Value = 142.8125
$ws.Range("K33").Value = 142.8125
$ws.Range("M33").Value = 86.1875

# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 4061.0789
$ws.Range("J40").Value = 3837.743
$ws.Range("L40").Value = 3837.743
$ws.Range("N40").Value = -4187.743

# row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 2501
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 166671540
$ws.Range("J64").Value = 500004380
$ws.Range("L64").Value = 500004380
$ws.Range("N64").Value = -500004876

# row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 166671540
$ws.Range("J67").Value = 500004380
$ws.Range("L67").Value = 500004380
$ws.Range("N67").Value = -500006096

# row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 1184.1428
$ws.Range("I92").Value = 1285.9
$ws.Range("K92").Value = 1285.9
$ws.Range("M92").Value = -37.90000000000009

# row 99 (Leve Item ID 19883)
$ws.Range("H99").Value = 316.2
$ws.Range("I99").Value = 351.75
$ws.Range("K99").Value = 1055.25
$ws.Range("M99").Value = 442.75

# row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 1648.95
$ws.Range("J112").Value = 1962.125
$ws.Range("L112").Value = 5886.375
$ws.Range("N112").Value = -8102.375

# row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1651.9231
$ws.Range("I129").Value = 559.75
$ws.Range("J129").Value = 3399.4
$ws.Range("K129").Value = 1679.25
$ws.Range("L129").Value = 10198.2
$ws.Range("M129").Value = 3320.75
$ws.Range("N129").Value = -20198.2

# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 8828.0625
$ws.Range("I132").Value = 8828.0625
$ws.Range("K132").Value = 26484.1875
$ws.Range("M132").Value = -23954.1875

# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2481.6667
$ws.Range("I138").Value = 1604.7142
$ws.Range("J138").Value = 2748.5652
$ws.Range("K138").Value = 4814.142599999999
$ws.Range("L138").Value = 8245.695599999999
$ws.Range("M138").Value = 325.8574000000008
$ws.Range("N138").Value = -18525.6956

# row 139 (Leve Item ID 42306)
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("ARM")
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 2290386.2
$ws.Range("I32").Value = 1124916.5
$ws.Range("K32").Value = 1124916.5
$ws.Range("M32").Value = -1124629.5

# row 133 (Leve Item ID 41857)
$ws.Range("H133").Value = 107999
$ws.Range("J133").Value = 107999
$ws.Range("L133").Value = 107999
$ws.Range("N133").Value = -113059

$ws = $wb.Worksheets.Item("BSM")
# row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 45462256
$ws.Range("I20").Value = 71438820
$ws.Range("K20").Value = 71438820
$ws.Range("M20").Value = -71438573

# row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 810.75
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 748
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 748
$ws.Range("M22").Value = -826
$ws.Range("N22").Value = -1094

# row 59 (Leve Item ID 43223)
$ws.Range("H59").Value = 115505.75
$ws.Range("J59").Value = 115505.75
$ws.Range("L59").Value = 115505.75
$ws.Range("N59").Value = -117199.75

$ws = $wb.Worksheets.Item("CRP")
# row 6 (Leve Item ID 2219)
$ws.Range("H6").Value = 25002500
$ws.Range("I6").Value = 25002500
$ws.Range("K6").Value = 25002500
$ws.Range("M6").Value = -25002387

# row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 216.25
$ws.Range("J7").Value = 400
$ws.Range("L7").Value = 400
$ws.Range("N7").Value = -626

# row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1020.5
$ws.Range("I16").Value = 1062.3334
$ws.Range("J16").Value = 895
$ws.Range("K16").Value = 1062.3334
$ws.Range("L16").Value = 895
$ws.Range("M16").Value = -775.3334
$ws.Range("N16").Value = -1469

# row 52 (Leve Item ID 43237)
$ws.Range("H52").Value = 52824
$ws.Range("J52").Value = 52824
$ws.Range("L52").Value = 52824
$ws.Range("N52").Value = -53412

# row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1020.5
$ws.Range("I113").Value = 1062.3334
$ws.Range("J113").Value = 895
$ws.Range("K113").Value = 1062.3334
$ws.Range("L113").Value = 895
$ws.Range("M113").Value = 1107.6666
$ws.Range("N113").Value = -5235

# row 138 (Leve Item ID 42302)
$ws.Range("H138").Value = 69998.91
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 69998.91
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 69998.91
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -80278.91

$ws = $wb.Worksheets.Item("CUL")
# row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 4314
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120

# row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 9151
$ws.Range("I137").Value = 1918.3334
$ws.Range("J137").Value = 20000
$ws.Range("K137").Value = 5755.0002
$ws.Range("L137").Value = 60000
$ws.Range("M137").Value = -655.0002000000004
$ws.Range("N137").Value = -70200

$ws = $wb.Worksheets.Item("GSM")
# row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2843.3333
$ws.Range("I113").Value = 2440
$ws.Range("K113").Value = 2440
$ws.Range("M113").Value = -270

# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2843.4814
$ws.Range("I132").Value = 2640.4707
$ws.Range("J132").Value = 3188.6
$ws.Range("K132").Value = 7921.4121
$ws.Range("L132").Value = 9565.799999999999
$ws.Range("M132").Value = -5391.4121
$ws.Range("N132").Value = -14625.8

$ws = $wb.Worksheets.Item("LTW")
# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 3950.5
$ws.Range("I22").Value = 3950.5
$ws.Range("M22").Value = -3655.5

# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 3950.5
$ws.Range("I27").Value = 3950.5
$ws.Range("K27").Value = 3950.5
$ws.Range("M27").Value = -3843.5

# row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 8534.125
$ws.Range("I46").Value = 1944
$ws.Range("J46").Value = 15124.25
$ws.Range("K46").Value = 1944
$ws.Range("L46").Value = 15124.25
$ws.Range("M46").Value = -1756
$ws.Range("N46").Value = -15500.25

# row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 708.6
$ws.Range("I55").Value = 459.8
$ws.Range("K55").Value = 459.8
$ws.Range("M55").Value = -286.8

# row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 6779.5356
$ws.Range("I122").Value = 5760
$ws.Range("J122").Value = 8614.700000000001
$ws.Range("K122").Value = 17280
$ws.Range("L122").Value = 25844.1
$ws.Range("M122").Value = -14830
$ws.Range("N122").Value = -30744.1

$ws = $wb.Worksheets.Item("WVR")
# row 14 (Leve Item ID 2658)
$ws.Range("H14").Value = 12500
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 578.5
$ws.Range("I113").Value = 578.5
$ws.Range("K113").Value = 1735.5
$ws.Range("M113").Value = 434.5
